# T1503_Contact_ContactDetails_AddEditDeleteActivity.xlsx
# Rewrites the "Contact" sample-data header/rows and replaces the
# "Activity" sheet's 2-column lookup table with the new 7-column
# Activity test-data grid (Type/Subject/IndustryGroup/ProductType/
# Description/MeetingNotes/ExtAttendee), matching the target diff.
#
# NOTE on ordering: the workbook's sharedStrings table is rebuilt on
# save by keeping still-referenced strings in their original relative
# order and appending newly-introduced strings in first-use order. To
# land on the exact shared-string layout of the target file, the cell
# writes below are intentionally sequenced: Contact!A1:C2 first, then
# the whole Activity sheet, then Contact!A3:C3 last.

$wb = $excel.ActiveWorkbook

$wsContact = $wb.Worksheets.Item("Contact")
$wsActivity = $wb.Worksheets.Item("Activity")

# --- Contact: header row + first data row -----------------------------
$wsContact.Range("A1").Value = "ContactName"
$wsContact.Range("B1").Value = "RelatedCompany"
$wsContact.Range("C1").Value = "Tab"

$wsContact.Range("A2").Value = "Activity Test External Contact"
$wsContact.Range("B2").Value = "ActivityCompany"
$wsContact.Range("C2").Value = "Activity"

# --- Activity: new 7-column grid ---------------------------------------
$wsActivity.Cells.Item(1, 1).Value = "Type"
$wsActivity.Cells.Item(1, 2).Value = "Subject"
$wsActivity.Cells.Item(1, 3).Value = "IndustryGroup"
$wsActivity.Cells.Item(1, 4).Value = "ProductType"
$wsActivity.Cells.Item(1, 5).Value = "Description"
$wsActivity.Cells.Item(1, 6).Value = "MeetingNotes"
$wsActivity.Cells.Item(1, 7).Value = "ExtAttendee"

$activityRows = @(
    @("Meeting", "Automated Test Subject Meeting", "BUS - Business Services", "Activist Advisory", "Automated Test Description Meeting", "Meeting Notes 1", "Test External"),
    @("Call",    "Automated Test Subject Call",    "BUS - Business Services", "Activist Advisory", "Automated Test Description Call",    "Meeting Notes 2", "Test External"),
    @("Email",   "Automated Test Subject Email",   "BUS - Business Services", "Activist Advisory", "Automated Test Description Email",   "Meeting Notes 3", "Test External"),
    @("Other",   "Automated Test Subject Other",   "BUS - Business Services", "Activist Advisory", "Automated Test Description Other",   "Meeting Notes 4", "Test External")
)

for ($r = 0; $r -lt $activityRows.Length; $r++) {
    $rowValues = $activityRows[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $wsActivity.Cells.Item($r + 2, $c + 1).Value = $rowValues[$c]
    }
}

# Header row: bold + centered (new cellXfs entry applyFont+applyAlignment)
$headerRow = $wsActivity.Range("A1:G1")
$headerRow.Font.Bold = $true
$headerRow.HorizontalAlignment = -4108   # xlCenter

# Column widths (best-effort; real Excel autosizes these to content)
$wsActivity.Columns.Item(1).AutoFit()
$wsActivity.Columns.Item(2).AutoFit()
$wsActivity.Columns.Item(3).AutoFit()
$wsActivity.Columns.Item(4).AutoFit()
$wsActivity.Columns.Item(5).AutoFit()
$wsActivity.Columns.Item(6).AutoFit()
$wsActivity.Columns.Item(7).AutoFit()

# --- Contact: final data row (written last so "Test Houlihan" is the ---
# --- last shared string introduced, matching the target ordering) -----
$wsContact.Range("A3").Value = "Test Houlihan"
$wsContact.Range("B3").Value = "StandardTestCompany"
$wsContact.Range("C3").Value = "Activity"

$wsContact.Columns.Item(1).AutoFit()
$wsContact.Columns.Item(2).AutoFit()

# --- Selections: touch Activity's selection first, then return to -----
# --- Contact so it ends up the active/selected tab, matching target ---
$wsActivity.Range("B14").Select()
$wsContact.Activate()
$wsContact.Range("C4").Select()
